# Update the 2025 row (row 7) of the faturamento_anual sheet with refreshed
# billing figures ("atualizei dados bibi e add").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2692931.66
$ws.Range("C7").Value = -39.39040486643525
$ws.Range("D7").Value = 2720
$ws.Range("E7").Value = 2720
$ws.Range("F7").Value = 990.0484044117647
$ws.Range("G7").Value = 5.53200093844215
